$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.627.78'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.59'
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  -2.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.82'
$ws.Range("E5").Value = '  -1.36%  '

$ws.Range("E6").Value = '  -1.95%  '

$ws.Range("E7").Value = '  -2.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3735'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07334'
$ws.Range("E9").Value = '  -1.76%  '

$ws.Range("E10").Value = '  -1.41%  '

$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.846.38'
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.709'
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.437'
$ws.Range("E14").Value = '  -1.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07128'
$ws.Range("E15").Value = '  -1.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.31'
$ws.Range("E16").Value = '  +4.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.015'
$ws.Range("E17").Value = '  -2.23%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008965'
$ws.Range("E18").Value = '  -1.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("E19").Value = '  -1.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.43'
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.632.04'
$ws.Range("E21").Value = '  -0.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.208'
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.07'
$ws.Range("E23").Value = '  -1.96%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.081.93'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("E25").Value = '  -0.52%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.64'
$ws.Range("E26").Value = '  -1.61%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.58'
$ws.Range("E27").Value = '  -1.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.153'
$ws.Range("E28").Value = '  +8.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.356'
$ws.Range("E29").Value = '  -0.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.49'
$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08937'
$ws.Range("E31").Value = '  -1.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.226'
$ws.Range("E32").Value = '  -0.14%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7745'
$ws.Range("E33").Value = '  -0.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.539'
$ws.Range("E34").Value = '  -1.31%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.877'
$ws.Range("E35").Value = '  -4.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.013'
$ws.Range("E36").Value = '  -2.07%  '

$ws.Range("E37").Value = '  -1.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05326'
$ws.Range("E38").Value = '  -0.60%  '

$ws.Range("E39").Value = '  -0.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.273'
$ws.Range("E40").Value = '  +5.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.921'
$ws.Range("E41").Value = '  +1.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5103'
$ws.Range("E42").Value = '  -2.14%  '

$ws.Range("E43").Value = '  -0.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.779'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.04'
$ws.Range("E45").Value = '  -1.50%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.65'
$ws.Range("E46").Value = '  -0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4733'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06479'
$ws.Range("E48").Value = '  -3.12%  '

$ws.Range("E49").Value = '  -2.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.688'
$ws.Range("E50").Value = '  -1.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.840'
$ws.Range("E51").Value = '  -4.04%  '
